$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$content = @"
h1: Justifying Imran Khan
p.note: Yesterday UN General Seceretary Antonio Guterres was visiting various places in Pakistan beside Imran Khan. I wanted to explain it religiously but it needs colors we all can see.
img.width-m-half: Antonio.png
p.caption.text-center: <a href=https://twitter.com/antonioguterres/status/1229730148964012034?s=20 target=_blank>Link to the tweet here</a>
p: Imran Khan, like Donald Trump, has always been the low cast, immature and spoiled kid in the world of politics. Donald Trump, I kept calling him, <b>he will fail</b>, did not fail. He made it all the way to the <b>President</b> of USA. Imran Khan is a similar story. Everyone called him a <b>failure</b> until we saw him rule 1 x province of Pakistan. But, the virus of PTI (Pakistan Tehreek Insaf) did not stop there. It seeped into other provinces. He is now the sitting <b>Prime Minister of Pakistan</b>.
h3: Why did he not fail?
p: He has not succeeded yet. But how did he move up there.
img.width-half: imrankhanbook.jpg
p: I was reading this biography of <a href=https://books.google.com.pk/books?id=5ptqDuf48RMC&dq=biography+of+imran+khan&hl=en&sa=X&ved=0ahUKEwjx1sKysNznAhVDgRoKHTLPCK0Q6AEILDAB target=_blank> Imran Khan by Christopher Sandford</a>; where he narrated an incident on the cricket ground. He called all the players for their opinion. Heard them. Acknowledged them. But went his own way. He always had a plan in his mind, <b>to execute</b>. Against all odds, suggestions and favors, <b>he walked on his own beaten road</b>. Invisible to his players, fellows and colleagues. He just walked with his head down to the glory of 1992 World Cup. 
p: In 1999, we had all the best cricket stars but we lacked Imran Khan in the field, beside the great players to unite them and bring them to glory. 
p: Imran Khan has been divorced twice, that means he did not even listen to his wives. <b>He did not stop</b>, is one reason that made him all things, he is today.
h3: Who all are failing / failed?
p: <b>Nawaz Sharif</b>, is on exile. Running for his life. Hiding in abroad. Wants to bring his daughter back, can’t convince anyone. He can’t return to Pakistan, he will go back to prison. 
p: <b>Zardari</b>, is the name of a land grabber mafia. Sitting in prison. Each passing day, he gets another blow on his reputation and freedom. Tired, failing and falling.
p: <b>Modi</b> is dejecting world leaders from entering India. They are denouncing him for all the brutalities and cruelties against minorities in India. Not just the Muslims, the educated lot is under attack. About to fall, dwindle and diminish.
p: <b>All the people</b> trying to bring Imran Khan down for his efforts, are failing everyday. 
p: <b>Opposition is greedy</b>, is another reason for Imran Khan’s so far success.
h3: Justifying Imran Khan
p: He is bringing delegations. Every day he is standing on the rostrum, speaking. Convincing us to have patience. Each day he is blown hard from all sides. He is re-affirmed, his actions are wrong. Each day people are unwilling to remove luxuries from their lives, cut down their finances, sacrifice on their benefits, label Imran Khan ‘Wrong’. Each passing moment he is dragged hard by us, we people who can not want to see a change.
p: He needs help. Bring your money, benefits, children to him. He talks of Madina, bring your life to him. He talks of Islam and its successes, what stops you from giving him the momentum he needs. He is going to win. Just few days ahead, he is going to convince the world of <b>the beauty</b> that lies in the foundations of our religion. The world might know, it might realize, it might get that sleek light that shined in Badr, Baghdad and Yarmouk. The world might feel, what it feels to live in harmony, rythem, abundance and love. They might see what the world of God looks like. Imran Khan might make this world the rock that Allah wants it to be. He might make millions proud of God. He might…
p: He is there to succeed. Say what ever you want to say. People came and went on. They tried to lower the Mujahideen in the past. They struggled to break the momentum. They failed. It rocked in the past. It will rock hard again. Bring it down. You can not.
h3: Conclusion
p: We be good, we succeed. We be greed, we fail. It is no rocket science but plain old simple math. Speak good, be kind and act sober. No power on earth can undo a good human being. No resistance can survive and beat this wall down. It goes up and that is the only way history has worked. Winston Churchill, against all odds and thugs, sits in a public transit, seeks opinion and acts. Might be drunk, but saves entire Kingdom.
p: It does not take to be a great scholar, philosopher or education, it takes one moment to become good. In a moment, we can undo all of ourselves. We can feel great and fight great. Just a moment, makes us a millionaire. It is not science. It is history speaking or it is <span class=lavendar>Quran</span> speaking.
quote: Those who believe fight in the cause of Allah, and those who disbelieve fight in the cause of Taghut. So fight against the allies of Satan. Indeed, <span class=lavendar>the plot of Satan has ever been weak</span>. <br> - Surah Al Nisa verse 76
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
"@

$ws.Cells.Item(50,1).Value2 = 49
$ws.Cells.Item(50,2).Value2 = 43879
$ws.Cells.Item(50,2).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(50,3).Value2 = "Surah Al Nisa, 68 - 87"
$ws.Cells.Item(50,6).Value2 = "Imran Khan, Modi, Nawaz Sharif, Zardari, The weaker plot, Pablo Escobar, Successful people"
$ws.Cells.Item(50,4).Value2 = $content
$ws.Cells.Item(50,5).Value2 = "Qasim Ali"

$ws.Rows.Item(50).RowHeight = 409.6

$win = $wb.Windows.Item(1)
$win.ScrollRow = 50
$win.ScrollColumn = 1
$ws.Range("D50").Select() | Out-Null
